# Generate Report for Handback
# Adds a new handback-report row for file c75b0bf8-5cfd-4387-a954-9685ab5a0fea
# to the Overview, zh-cn and de-de sheets (row 4 on each sheet).

$wb = $excel.ActiveWorkbook

$uuid = "c75b0bf8-5cfd-4387-a954-9685ab5a0fea"
$xlfHash = "15140a5e94629eb678e55ae3f0fc69efb766c8b8"

$mdCommit          = "215718091fdd18325e89cb5d192c32ec8caa6f8c"
$handoffCommitZh   = "57afe40a64ca63d158e638f443432909a8ddd9b9"
$zhMdCommit        = "2a74d1819f01d4bdac78f6823abd7663249e5f10"
$handbackCommitZh  = "13698b94b1c29a34ce48908c61cd31699aa93c52"
$handoffCommitDe   = "5396cb07b69ae0afb8ce25e22daa21e09bc4b6d9"
$deMdCommit        = "b0a44bd814d96d0c45c79296627561a7dc245f2a"
$handbackCommitDe  = "c0f12f71ce56f6bfdcccfd765158a08eb786af80"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$uuid.md"

$handoffXlfZh  = "$uuid.$xlfHash.zh-cn.xlf"
$handoffXlfDe  = "$uuid.$xlfHash.de-de.xlf"

$handoffUrlZh  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffCommitZh/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$handoffXlfZh"
$zhMdUrl       = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$zhMdCommit/e2e/$uuid.md"
$handbackUrlZh = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$handbackCommitZh/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$handoffXlfZh"

$handoffUrlDe  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffCommitDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$handoffXlfDe"
$deMdUrl       = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$deMdCommit/e2e/$uuid.md"
$handbackUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$handbackCommitDe/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$handoffXlfDe"

$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrl, "", "", "$uuid.md") | Out-Null
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdUrl, "", "", "$uuid.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), $mdUrl, "", "", ".md") | Out-Null
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $handoffUrlZh, "", "", $handoffXlfZh) | Out-Null
$wsZh.Range("E4").Value = "2016-03-09 10:39:16"
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), $zhMdUrl, "", "", "$uuid.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G4"), $handbackUrlZh, "", "", $handoffXlfZh) | Out-Null
$wsZh.Range("H4").Value = "2016-03-09 10:40:10"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdUrl, "", "", "$uuid.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), $mdUrl, "", "", ".md") | Out-Null
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $handoffUrlDe, "", "", $handoffXlfDe) | Out-Null
$wsDe.Range("E4").Value = "2016-03-09 10:39:27"
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), $deMdUrl, "", "", "$uuid.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G4"), $handbackUrlDe, "", "", $handoffXlfDe) | Out-Null
$wsDe.Range("H4").Value = "2016-03-09 10:40:29"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = "Include"
